# Truth-table workbook touch-up: a few cell corrections in the row-4
# example, a couple of numeric fills further down, unmerging the
# label cells in columns E:G (and C:G for the two immediate rows),
# and moving the live selection to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 example: rd/rt/X got swapped around ---
$ws.Range("D4").Value = "rt"
$ws.Range("I4").Value = "rt"
$ws.Range("J4").Value = "X"

# --- Fill in a couple of previously-blank/incorrect binary codes ---
$ws.Range("M27").Value = 1100
$ws.Range("G30").Value = 100000

# --- Remove all the merged "——Immediate——" / label cells ---
$ws.Range("E4:G4").UnMerge()
$ws.Range("E12:G12").UnMerge()
$ws.Range("E13:G13").UnMerge()
$ws.Range("E14:G14").UnMerge()
$ws.Range("E15:G15").UnMerge()
$ws.Range("C16:G16").UnMerge()
$ws.Range("C17:G17").UnMerge()
$ws.Range("E25:G25").UnMerge()
$ws.Range("E26:G26").UnMerge()
$ws.Range("E27:G27").UnMerge()
$ws.Range("E28:G28").UnMerge()

# --- Move the selection ---
$ws.Range("A6").Select()
$ws.Range("E5").Select()
